# "ppt is sub-task done"
#
# 1) Update every "Update automatically" date field in the deck
#    (slides, slide layouts, slide master, notes master) from the
#    7 March 2022 / 3-7-2022 capture date to 13 March 2022 / 3-13-2022.
# 2) Slide 1: merge the "More " / "language support?" runs into a
#    single run "More language support?".
# 3) Slide 2: merge the " server: " / "I/O to file (store data)" runs
#    into a single run " server: I/O to file (store data)" (leaving the
#    preceding "nodeJs" run untouched).

$p = $ppt.ActivePresentation

function Set-DateShapeText {
    param($shapes, [string]$marker, [string]$newText)
    for ($j = 1; $j -le $shapes.Count; $j++) {
        $shp = $shapes.Item($j)
        if ($shp.HasTextFrame) {
            $txt = $shp.TextFrame.TextRange.Text
            if ($txt -like "*$marker*") {
                $shp.TextFrame.TextRange.Text = $newText
                return $true
            }
        }
    }
    return $false
}

# --- 1a) Slides: "7 March 2022" -> "13 March 2022" ---------------------
for ($k = 1; $k -le $p.Slides.Count; $k++) {
    $s = $p.Slides.Item($k)
    Set-DateShapeText $s.Shapes "7 March 2022" "13 March 2022" | Out-Null
}

# --- 1b) Slide layouts: "7 March 2022" -> "13 March 2022" --------------
$layouts = $p.SlideMaster.CustomLayouts
for ($k = 1; $k -le $layouts.Count; $k++) {
    $lyt = $layouts.Item($k)
    Set-DateShapeText $lyt.Shapes "7 March 2022" "13 March 2022" | Out-Null
}

# --- 1c) Slide master: "7 March 2022" -> "13 March 2022" ---------------
Set-DateShapeText $p.SlideMaster.Shapes "7 March 2022" "13 March 2022" | Out-Null

# --- 1d) Notes master: "3/7/2022" -> "3/13/2022" ------------------------
# (Direct Shapes(..).TextFrame edits don't stick on the notes master in
# this host; the HeadersFooters.DateAndTime property is the reliable path.)
$p.NotesMaster.HeadersFooters.DateAndTime.Text = "3/13/2022"

# --- 2) Slide 1: merge "More " + "language support?" -------------------
$slide1 = $p.Slides.Item(1)
for ($j = 1; $j -le $slide1.Shapes.Count; $j++) {
    $shp = $slide1.Shapes.Item($j)
    if ($shp.HasTextFrame) {
        $tr = $shp.TextFrame.TextRange
        if ($tr.Text -like "*More language support?*") {
            for ($i = 1; $i -le $tr.Paragraphs(0,0).Count + 50; $i++) {
                $para = $tr.Paragraphs($i, 1)
                if ($para.Length -eq 0) { break }
                if ($para.Text -eq "More language support?") {
                    # Force a real text change so the run-merge is applied
                    # (identical-text assignments are treated as no-ops),
                    # then set the final text using a freshly-fetched range.
                    $para.Text = "More language support?__tmp__"
                    $tr2 = $shp.TextFrame.TextRange
                    $para2 = $tr2.Paragraphs($i, 1)
                    $para2.Text = "More language support?"
                    break
                }
            }
            break
        }
    }
}

# --- 3) Slide 2: merge " server: " + "I/O to file (store data)" --------
$slide2 = $p.Slides.Item(2)
for ($j = 1; $j -le $slide2.Shapes.Count; $j++) {
    $shp = $slide2.Shapes.Item($j)
    if ($shp.HasTextFrame) {
        $tr = $shp.TextFrame.TextRange
        if ($tr.Text -like "*nodeJs server: I/O to file (store data)*") {
            for ($i = 1; $i -le $tr.Paragraphs(0,0).Count + 50; $i++) {
                $para = $tr.Paragraphs($i, 1)
                if ($para.Length -eq 0) { break }
                if ($para.Text -eq "nodeJs server: I/O to file (store data)") {
                    $start = $para.Start
                    $len = $para.Length
                    $prefixLen = 6  # "nodeJs"
                    $sub = $tr.Characters($start + $prefixLen, $len - $prefixLen)
                    $sub.Text = " server: I/O to file (store data)__tmp__"

                    $tr2 = $shp.TextFrame.TextRange
                    $para2 = $tr2.Paragraphs($i, 1)
                    $sub2 = $tr2.Characters($para2.Start + $prefixLen, $para2.Length - $prefixLen)
                    $sub2.Text = " server: I/O to file (store data)"
                    break
                }
            }
            break
        }
    }
}
